# Fixing Coins and Gems names
# - A7: clarify the "recruitingDuration" note -> "recruitingLimitDuration"
# - Reward rows that pay out coins/gems: [type] column (H) now uses the
#   short codes "sc" / "hc", and the previously-blank [rewardSku] column (J)
#   is filled in with the human-readable "Coins" / "Gems" label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# SETTINGS note text fix
$ws.Range("A7").Value = "recruitingLimitDuration: Days when players can start the xpromo, starting from [startDate] (this is a date also)"

# reward_hd_hse_1b (row 16): coins -> sc, rewardSku Coins
$ws.Range("H16").Value = "sc"
$ws.Range("J16").Value = "Coins"

# reward_hd_hse_2b (row 17): gems -> hc, rewardSku Gems
$ws.Range("H17").Value = "hc"
$ws.Range("J17").Value = "Gems"

# reward_hd_hse_3b (row 18): gems -> hc, rewardSku Gems
$ws.Range("H18").Value = "hc"
$ws.Range("J18").Value = "Gems"

# reward_hse_hd_1b (row 22): coins -> sc, rewardSku Coins
$ws.Range("H22").Value = "sc"
$ws.Range("J22").Value = "Coins"

# reward_hse_hd_2b (row 23): gems -> hc, rewardSku Gems
$ws.Range("H23").Value = "hc"
$ws.Range("J23").Value = "Gems"

# Move the active selection to C9 (matches the author's cursor position in
# the committed workbook)
$ws.Activate()
$ws.Range("C9").Select()
